# Update countries & provincias Spain
# Updates the COVID data snapshot in the "Pais" sheet: refreshes the
# "updated at" timestamp, refreshes a handful of per-country stat rows,
# and re-labels a few rows whose ranking order changed (their A-column
# country name swaps with the neighboring row, carrying along their own
# stats) as a result of the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 14:57"

# --- Refreshed per-country statistics (no reordering) ------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 2553771
$ws.Range("C4").Value = 815
$ws.Range("E4").Value = 1357354

# Row 18: Arabia Saudita
$ws.Range("B18").Value = 178504
$ws.Range("C18").Value = 3927
$ws.Range("D18").Value = 122128
$ws.Range("E18").Value = 54865
$ws.Range("G18").Value = 37
$ws.Range("H18").Value = 1511

# Row 33: Paises Bajos
$ws.Range("B33").Value = 50074
$ws.Range("C33").Value = 69
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 6105

# Row 50: Barein
$ws.Range("E50").Value = 5590
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 78

# Row 77: Consejo Danes para los Refugiados
$ws.Range("B77").Value = 6690
$ws.Range("C77").Value = 138
$ws.Range("D77").Value = 937
$ws.Range("E77").Value = 5600
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 153

# Row 78: Senegal
$ws.Range("B78").Value = 6459
$ws.Range("C78").Value = 105
$ws.Range("D78").Value = 4255
$ws.Range("E78").Value = 2102
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 102

# Row 79: Republica de Macedonia
$ws.Range("B79").Value = 5906
$ws.Range("C79").Value = 148
$ws.Range("D79").Value = 2236
$ws.Range("E79").Value = 3393
$ws.Range("G79").Value = 9
$ws.Range("H79").Value = 277

# Row 101: Croacia
$ws.Range("B101").Value = 2624
$ws.Range("C101").Value = 85
$ws.Range("D101").Value = 2152
$ws.Range("E101").Value = 365

# --- Rows whose relative order/labels changed --------------------------
# Laos / Santa Lucia swap position (tied totals, new tie-break order)
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Laos"

# Dominica / Fiyi swap position (tied totals, new tie-break order)
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A204").Value = "Dominica"

# Islas Malvinas / Groenlandia swap position (tied totals)
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

# Montserrat / Seychelles swap position, bringing their own stats with them
$ws.Range("A212").Value = "Seychelles"
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 10
$ws.Range("H213").Value = 1
